$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.440.93"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.913.83"
$ws.Range("E3").Value = "  +3.68%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.44"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.83"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.559"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.630"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.96"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.89"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.80"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.369.02"
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("E16").Value = "  +6.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.927.80"
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.434.57"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.64"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  +3.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.15"
$ws.Range("E21").Value = "  +3.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.88"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.17"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.78"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.76"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("E27").Value = "  +3.45%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.88"
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.34"
$ws.Range("E31").Value = "  +12.28%  "
$ws.Range("E32").Value = "  +7.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0984"
$ws.Range("E33").Value = "  +11.40%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.46"
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  +5.48%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("E40").Value = "  +2.68%  "
$ws.Range("E41").Value = "  +13.79%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.70"
$ws.Range("E42").Value = "  +7.63%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.117"
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("E44").Value = "  +9.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "120.94"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("E47").Value = "  +3.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.200.94"
$ws.Range("E48").Value = "  +4.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.266"
$ws.Range("E49").Value = "  +24.26%  "
$ws.Range("B50").Value = "BEAM"
$ws.Range("C50").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0340"
$ws.Range("E50").Value = "  +12.42%  "
$ws.Range("B51").Value = "SEI"
$ws.Range("C51").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.976"
$ws.Range("E51").Value = "  +3.11%  "

Write-Output "Applied 93 cell updates"
